# Fruta / hortaliza, semanal
# Insert a new weekly price-report row at row 25 of the "Poroto granado" sheet.
# This shifts the existing rows 25-76 down to 26-77 (dimension A1:R76 -> A1:R77)
# and populates the newly inserted row with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 25, pushing everything else down.
$ws.Rows.Item(25).Insert()

# Fill the new row 25 with the new weekly observation.
$ws.Cells.Item(25, 1).Value = 10
$ws.Cells.Item(25, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(25, 3).Value = "La Araucanía"
$ws.Cells.Item(25, 4).Value = 44614
$ws.Cells.Item(25, 5).Value = 9
$ws.Cells.Item(25, 6).Value = 100112030
$ws.Cells.Item(25, 7).Value = "Poroto granado"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 40
$ws.Cells.Item(25, 11).Value = 25000
$ws.Cells.Item(25, 12).Value = 25000
$ws.Cells.Item(25, 13).Value = 25000
$ws.Cells.Item(25, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(25, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(25, 16).Value = 1000
$ws.Cells.Item(25, 17).Value = 25
$ws.Cells.Item(25, 18).Value = "Hortaliza"
